$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A20:E20").Style = "Normal"

$ws.Range("A20").Value = "release/8.0.17"
$ws.Range("B20").Value = "X"
$ws.Range("C20").Value = "X"
$ws.Range("D20").Value = "X"
$ws.Range("E20").Value = "X"
